# Update "Handed back"/report generation timestamps produced when the
# handback report was regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-29 13:06:21"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-29 13:06:16"
$wsZhCn.Range("K2").Value = "2016-08-29 13:06:34"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-29 13:06:21"
$wsDeDe.Range("K2").Value = "2016-08-29 13:06:41"
